$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = '28.128.57'
$ws.Range("E2").Value = '  -3.34%  '
$ws.Range("D3").Value = '1.926.27'
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("E4").Value = '  -0.63%  '
Set-TextValue $ws.Range("D5") '329.46'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").Value = '  -0.66%  '
Set-TextValue $ws.Range("D7") '0.4728'
$ws.Range("E7").Value = '  -4.98%  '
Set-TextValue $ws.Range("D8") '0.4061'
$ws.Range("E8").Value = '  -3.64%  '
Set-TextValue $ws.Range("D9") '52.91'
$ws.Range("E9").Value = '  -0.49%  '
Set-TextValue $ws.Range("D10") '0.08434'
$ws.Range("E10").Value = '  -9.04%  '
Set-TextValue $ws.Range("D11") '1.048'
$ws.Range("E11").Value = '  -4.79%  '
Set-TextValue $ws.Range("D12") '22.30'
$ws.Range("E12").Value = '  -2.61%  '
$ws.Range("D13").Value = '1.937.14'
$ws.Range("E13").Value = '  -1.63%  '
Set-TextValue $ws.Range("D14") '7.515'
$ws.Range("E14").Value = '  -5.00%  '
Set-TextValue $ws.Range("D15") '6.102'
$ws.Range("E15").Value = '  -5.59%  '
$ws.Range("E16").Value = '  -0.87%  '
Set-TextValue $ws.Range("D17") '90.52'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("E18").Value = '  -3.79%  '
Set-TextValue $ws.Range("D19") '0.06576'
$ws.Range("E19").Value = '  -2.24%  '
Set-TextValue $ws.Range("D20") '18.09'
$ws.Range("E20").Value = '  -6.20%  '
Set-TextValue $ws.Range("D21") '1.002'
$ws.Range("E21").Value = '  -0.63%  '
Set-TextValue $ws.Range("D22") '5.756'
$ws.Range("E22").Value = '  -3.67%  '
$ws.Range("D23").Value = '28.142.75'
$ws.Range("E23").Value = '  -3.39%  '
Set-TextValue $ws.Range("D24") '11.41'
$ws.Range("E24").Value = '  -4.75%  '
Set-TextValue $ws.Range("D25") '2.285'
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("D26").Value = '2.162.81'
$ws.Range("E26").Value = '  -2.25%  '
Set-TextValue $ws.Range("D27") '154.03'
$ws.Range("E27").Value = '  -0.83%  '
Set-TextValue $ws.Range("D28") '20.12'
$ws.Range("E28").Value = '  -3.06%  '
Set-TextValue $ws.Range("D29") '2.155'
$ws.Range("E29").Value = '  -4.99%  '
Set-TextValue $ws.Range("D30") '5.716'
$ws.Range("E30").Value = '  -10.23%  '
Set-TextValue $ws.Range("D31") '123.71'
$ws.Range("E31").Value = '  -2.57%  '
Set-TextValue $ws.Range("D32") '0.9788'
$ws.Range("E32").Value = '  -6.76%  '
Set-TextValue $ws.Range("D33") '0.09602'
$ws.Range("E33").Value = '  -2.56%  '
Set-TextValue $ws.Range("D34") '1.450'
$ws.Range("E34").Value = '  -4.53%  '
Set-TextValue $ws.Range("D35") '5.562'
$ws.Range("E35").Value = '  -4.47%  '
$ws.Range("E36").Value = '  -2.60%  '
Set-TextValue $ws.Range("D37") '9.064'
$ws.Range("E37").Value = '  +0.04%  '
Set-TextValue $ws.Range("D38") '0.02318'
$ws.Range("E38").Value = '  -4.58%  '
Set-TextValue $ws.Range("D39") '0.06184'
$ws.Range("E39").Value = '  -3.77%  '
Set-TextValue $ws.Range("D40") '1.239'
$ws.Range("E40").Value = '  -6.40%  '
Set-TextValue $ws.Range("D41") '0.6177'
$ws.Range("E41").Value = '  -4.75%  '
Set-TextValue $ws.Range("D42") '11.06'
$ws.Range("E42").Value = '  -4.01%  '
Set-TextValue $ws.Range("D43") '1.001'
$ws.Range("E43").Value = '  -0.60%  '
Set-TextValue $ws.Range("D44") '0.1907'
$ws.Range("E44").Value = '  -4.78%  '
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D45") '1.295'
$ws.Range("E45").Value = '  -4.88%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range("D46") '0.5888'
$ws.Range("E46").Value = '  -5.39%  '
Set-TextValue $ws.Range("D47") '12.85'
$ws.Range("E47").Value = '  -3.44%  '
Set-TextValue $ws.Range("D48") '2.036'
$ws.Range("E48").Value = '  -6.99%  '
$ws.Range("E49").Value = '  -0.42%  '
Set-TextValue $ws.Range("D50") '0.06822'
$ws.Range("E50").Value = '  -2.04%  '
Set-TextValue $ws.Range("D51") '110.17'
$ws.Range("E51").Value = '  -2.73%  '
